$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.544.89"
$ws.Cells.Item(2, 5).Value = "  -0.38%  "
$ws.Cells.Item(3, 4).Value = "3.141.38"
$ws.Cells.Item(3, 5).Value = "  -0.56%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "575.10"
$ws.Cells.Item(5, 5).Value = "  +0.30%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "148.37"
$ws.Cells.Item(6, 5).Value = "  -1.87%  "
$ws.Cells.Item(7, 5).Value = "  +0.07%  "
$ws.Cells.Item(8, 4).Value = "3.141.98"
$ws.Cells.Item(8, 5).Value = "  -0.43%  "
$ws.Cells.Item(9, 5).Value = "  -0.68%  "
$ws.Cells.Item(10, 5).Value = "  -3.01%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "6.10"
$ws.Cells.Item(11, 5).Value = "  -1.49%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.498"
$ws.Cells.Item(12, 5).Value = "  -1.23%  "
$ws.Cells.Item(13, 5).Value = "  +1.58%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "37.02"
$ws.Cells.Item(14, 5).Value = "  -2.46%  "
$ws.Cells.Item(15, 4).Value = "3.659.32"
$ws.Cells.Item(15, 5).Value = "  -0.31%  "
$ws.Cells.Item(16, 4).Value = "64.670.99"
$ws.Cells.Item(16, 5).Value = "  -0.32%  "
$ws.Cells.Item(17, 4).Value = "3.141.09"
$ws.Cells.Item(17, 5).Value = "  -0.44%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "7.10"
$ws.Cells.Item(18, 5).Value = "  -1.98%  "
$ws.Cells.Item(19, 5).Value = "  +0.22%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "502.29"
$ws.Cells.Item(20, 5).Value = "  -2.40%  "
$ws.Cells.Item(21, 5).Value = "  -1.24%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.710"
$ws.Cells.Item(22, 5).Value = "  -3.80%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "15.13"
$ws.Cells.Item(23, 5).Value = "  -0.09%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "7.68"
$ws.Cells.Item(24, 5).Value = "  -2.32%  "
$ws.Cells.Item(25, 5).Value = "  -1.41%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  +0.00%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.89"
$ws.Cells.Item(27, 5).Value = "  -1.49%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.84"
$ws.Cells.Item(28, 5).Value = "  +1.23%  "
$ws.Cells.Item(29, 5).Value = "  -1.22%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.78"
$ws.Cells.Item(30, 5).Value = "  +4.62%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "27.45"
$ws.Cells.Item(31, 5).Value = "  -2.10%  "
$ws.Cells.Item(32, 5).Value = "  +0.06%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.20"
$ws.Cells.Item(33, 5).Value = "  +0.82%  "
$ws.Cells.Item(34, 5).Value = "  +0.25%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.43"
$ws.Cells.Item(35, 5).Value = "  -2.40%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "54.47"
$ws.Cells.Item(36, 5).Value = "  -2.23%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0887"
$ws.Cells.Item(37, 5).Value = "  +2.54%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "476.25"
$ws.Cells.Item(38, 5).Value = "  -2.09%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0413"
$ws.Cells.Item(39, 5).Value = "  -2.39%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.91"
$ws.Cells.Item(40, 5).Value = "  -3.77%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "8.67"
$ws.Cells.Item(41, 5).Value = "  +0.09%  "
$ws.Cells.Item(42, 4).Value = "2.995.99"
$ws.Cells.Item(42, 5).Value = "  -3.92%  "
$ws.Cells.Item(43, 5).Value = "  -4.77%  "
$ws.Cells.Item(44, 5).Value = "  -5.66%  "
$ws.Cells.Item(45, 5).Value = "  -2.66%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "27.96"
$ws.Cells.Item(46, 5).Value = "  -4.31%  "
$ws.Cells.Item(47, 4).Value = "0.0₃0577"
$ws.Cells.Item(47, 5).Value = "  -0.72%  "
$ws.Cells.Item(49, 5).Value = "  -1.96%  "
$ws.Cells.Item(50, 5).Value = "  -3.13%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "32.83"
$ws.Cells.Item(51, 5).Value = "  +3.79%  "
